$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-12 Saturday" "2023-08-13 Sunday"

Replace-Text "54÷6=" "50÷2="
Replace-Text "90÷8=" "68÷3="
Replace-Text "64÷7=" "26÷8="
Replace-Text "94÷2=" "73÷6="
Replace-Text "22÷3=" "27÷8="
Replace-Text "39÷2=" "15÷7="
Replace-Text "52÷3=" "84÷7="
Replace-Text "55÷3=" "38÷9="
Replace-Text "19÷7=" "98÷9="
Replace-Text "88÷8=" "31÷3="
Replace-Text "59÷7=" "74÷2="
Replace-Text "55÷7=" "92÷6="
Replace-Text "13÷8=" "29÷6="
Replace-Text "72÷2=" "33÷3="
Replace-Text "27÷5=" "20÷9="
Replace-Text "37÷5=" "38÷5="
Replace-Text "84÷5=" "91÷4="
Replace-Text "25÷7=" "48÷8="
Replace-Text "24÷5=" "22÷9="
Replace-Text "38÷3=" "55÷6="
Replace-Text "78÷8=" "56÷8="
Replace-Text "67÷2=" "50÷4="
Replace-Text "77÷5=" "72÷7="
Replace-Text "26÷7=" "21÷9="
Replace-Text "68÷2=" "28÷9="
